# Timetable update to V1.5
# Moves/swaps a handful of class entries between slots on Monday, Tuesday and
# Wednesday, and bumps the Wednesday sheet's zoom level.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Monday
# ---------------------------------------------------------------------------
$wsMon = $wb.Worksheets.Item("Monday")

# H11 ("ITC-G ... M. Shahzad") moves up into H5; H11 is cleared.
$monH11 = $wsMon.Range("H11").Value2
$wsMon.Range("H5").Value2 = $monH11
$wsMon.Range("H11").Value2 = $null

# G24 ("COAL Lab-C  Zain") and G30 ("DB-Lab C  Basit ali") swap places.
$monG24 = $wsMon.Range("G24").Value2
$monG30 = $wsMon.Range("G30").Value2
$wsMon.Range("G24").Value2 = $monG30
$wsMon.Range("G30").Value2 = $monG24

# G27 ("DB-Lab D  Ammara") moves up into G26; G27 is cleared.
$monG27 = $wsMon.Range("G27").Value2
$wsMon.Range("G26").Value2 = $monG27
$wsMon.Range("G27").Value2 = $null

# ---------------------------------------------------------------------------
# Tuesday
# ---------------------------------------------------------------------------
$wsTue = $wb.Worksheets.Item("Tuesday")

# Capture all source values before any writes so nothing gets clobbered.
$tueB5  = $wsTue.Range("B5").Value2   # Discrete-E  Nouman Durrani
$tueD5  = $wsTue.Range("D5").Value2   # Discrete-F  Nouman Durrani
$tueF5  = $wsTue.Range("F5").Value2   # Discrete-C  Nouman Durrani
$tueD14 = $wsTue.Range("D14").Value2  # NP-GR1  Shahbaz
$tueI7  = $wsTue.Range("I7").Value2   # COAL-E  Danish
$tueB24 = $wsTue.Range("B24").Value2  # COAL Lab-A  Zain
$tueB27 = $wsTue.Range("B27").Value2  # DB-Lab F  Ammara

# Row 5 entries redistribute to rows 7, 12 and 14.
$wsTue.Range("B5").Value2 = $null
$wsTue.Range("D5").Value2 = $null
$wsTue.Range("F5").Value2 = $null

$wsTue.Range("D7").Value2  = $tueD14
$wsTue.Range("F12").Value2 = $tueF5
$wsTue.Range("B14").Value2 = $tueB5
$wsTue.Range("D14").Value2 = $tueD5

# I7 moves down into I9; I7 is cleared.
$wsTue.Range("I9").Value2 = $tueI7
$wsTue.Range("I7").Value2 = $null

# B24 and B27 swap places.
$wsTue.Range("B24").Value2 = $tueB27
$wsTue.Range("B27").Value2 = $tueB24

# ---------------------------------------------------------------------------
# Wednesday
# ---------------------------------------------------------------------------
$wsWed = $wb.Worksheets.Item("Wednesday")

# C28 ("DB-Lab G  Ammara") moves up into C24; C28 is cleared.
$wedC28 = $wsWed.Range("C28").Value2
$wsWed.Range("C24").Value2 = $wedC28
$wsWed.Range("C28").Value2 = $null

# Bump the on-screen zoom level for the Wednesday sheet from 60% to 62%,
# without disturbing which sheet/tab is actually active in the workbook.
$originalActiveSheetName = $wb.ActiveSheet.Name
$wsWed.Activate()
$excel.ActiveWindow.Zoom = 62
$wb.Worksheets.Item($originalActiveSheetName).Activate()
